$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.291.61'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '2.306.37'
$ws.Range('E3').Value = '  -1.14%  '
$ws.Range('E4').Value = '  -0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.85'
$ws.Range('E5').Value = '  +2.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.57'
$ws.Range('E6').Value = '  -2.88%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.633'
$ws.Range('E7').Value = '  +0.39%  '
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.611'
$ws.Range('E9').Value = '  +0.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.77'
$ws.Range('E10').Value = '  -0.69%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0910'
$ws.Range('E11').Value = '  -0.96%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.36'
$ws.Range('E12').Value = '  -0.63%  '
$ws.Range('E13').Value = '  +0.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.969'
$ws.Range('E14').Value = '  -1.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.36'
$ws.Range('E15').Value = '  -2.14%  '
$ws.Range('D16').Value = '2.652.23'
$ws.Range('E16').Value = '  -1.71%  '
$ws.Range('D17').Value = '2.300.16'
$ws.Range('E17').Value = '  -2.25%  '
$ws.Range('D18').Value = '42.273.41'
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.47'
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('E20').Value = '  +0.56%  '
$ws.Range('E21').Value = '  +4.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.60'
$ws.Range('E22').Value = '  -2.97%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '281.50'
$ws.Range('E23').Value = '  +6.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.10'
$ws.Range('E24').Value = '  +18.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.28'
$ws.Range('E25').Value = '  -0.69%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.43%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.91'
$ws.Range('E27').Value = '  -2.35%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.39'
$ws.Range('E28').Value = '  +6.07%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '23.05'
$ws.Range('E29').Value = '  -0.92%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.25'
$ws.Range('E30').Value = '  +1.08%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '165.11'
$ws.Range('E31').Value = '  -0.56%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0876'
$ws.Range('E32').Value = '  -2.39%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.92'
$ws.Range('E33').Value = '  -1.08%  '
$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.136'
$ws.Range('E34').Value = '  +4.64%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.68'
$ws.Range('E35').Value = '  -8.25%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.115'
$ws.Range('E36').Value = '  -4.13%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.64'
$ws.Range('E37').Value = '  +1.27%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0350'
$ws.Range('E38').Value = '  -1.52%  '
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.76'
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.77'
$ws.Range('E40').Value = '  +5.18%  '
$ws.Range('B41').Value = 'BitcoinSV'
$ws.Range('C41').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '100.42'
$ws.Range('E41').Value = '  -3.54%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.47'
$ws.Range('E42').Value = '  -0.14%  '
$ws.Range('B43').Value = 'MultiversX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '69.35'
$ws.Range('E43').Value = '  -2.13%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.227'
$ws.Range('E44').Value = '  -3.51%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('B46').Value = 'Celestia'
$ws.Range('C46').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.05'
$ws.Range('E46').Value = '  -1.17%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '112.23'
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('B48').Value = 'ordi'
$ws.Range('C48').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '77.51'
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.92'
$ws.Range('E49').Value = '  -1.13%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.30'
$ws.Range('E50').Value = '  -3.02%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '1.609.20'
$ws.Range('E51').Value = '  +3.62%  '
